# Update of all scripts and data
#
# Row 69 (Gear "2-RAP", SpecCode "SCYOCAN", MatStage "MEDSE-1") is removed.
# All the following rows (70-75) shift up by one, turning the previous
# 7-row block (69-75) into a 6-row block (69-74); the sheet's used range
# shrinks from A1:Q75 to A1:Q74.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(69).Delete()
